$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.67856228351593
$ws.Range("B1").Value = 1.870242118835449
$ws.Range("C1").Value = 5.136736392974854
$ws.Range("D1").Value = 1.344742774963379
$ws.Range("E1").Value = 0.7432289719581604
